# Generate Report for Archive
#
# 1) The status text "Ready for handoff" becomes "In Translation" everywhere
#    it appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 - all backed by
#    the same shared-string entry).
# 2) Because the new status text is shorter, Excel's column auto-fit
#    shrinks the affected "Status" columns (Overview columns E/F, and
#    column C on the zh-cn / de-de sheets) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status text -------------------------------------------------
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# --- Re-fit the now-narrower "Status" columns -------------------------------
# (offset compensates for this host's fixed column-width padding so the
# resulting stored width lands as close as possible to Excel's real
# auto-fit result for the new, shorter text)
$targetWidth = 13.4101845877511
$padding = 5 / 6
$colInput = $targetWidth - $padding

$wsOverview.Columns.Item(5).ColumnWidth = $colInput  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $colInput  # column F
$wsZhCn.Columns.Item(3).ColumnWidth     = $colInput  # column C
$wsDeDe.Columns.Item(3).ColumnWidth     = $colInput  # column C
